$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---
# Force text type via NumberFormat "@" so numeric-looking strings (e.g. "23.15")
# are not silently coerced to numbers, then restore the default "Normal" style
# so no stray number-format style is left behind on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.188.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.471.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.471.55"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.067.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.471.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.203.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.566"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.618.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "32.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "169.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.510.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0765"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.799"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.602.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.15"
$ws.Range("D50").Style = "Normal"

# --- Volume(1h) column (E) updates ---
# These are percentage strings with two leading/trailing spaces; they are
# already stored as text so a plain Value assignment is sufficient.
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("E12").Value = "  +4.59%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("E20").Value = "  +3.64%  "
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  -5.09%  "
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("E38").Value = "  +8.62%  "
$ws.Range("E39").Value = "  +11.98%  "
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("E44").Value = "  +3.52%  "
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("E51").Value = "  +2.77%  "
